$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Priority Life Stages"),
# shifting Priority Life Stages -> G, Limiting Factor -> H, Action Categories -> I
$ws.Columns("F:F").Insert()

# Set the new header for column F
$ws.Range("F1").Value = "Action"

# Populate the new Action column for each data row
$ws.Range("F2").Value = 'Restore Fish Passage'
$ws.Range("F3").Value = ' Restore Reach Function'
$ws.Range("F4").Value = ' Restore Reach Function'
$ws.Range("F5").Value = ' Restore Reach Function'
$ws.Range("F6").Value = ' Restore Reach Function'
$ws.Range("F7").Value = ' Restore Reach Function'
$ws.Range("F8").Value = ' Restore Reach Function'
$ws.Range("F9").Value = ' Restore Reach Function'
$ws.Range("F10").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F11").Value = ' Restore Reach Function'
$ws.Range("F12").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F13").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F14").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F15").Value = ' Restore Reach Function'
$ws.Range("F16").Value = ' Restore Reach Function'
$ws.Range("F17").Value = ' Restore Reach Function'
$ws.Range("F18").Value = ' Restore Reach Function'
$ws.Range("F19").Value = ' Restore Reach Function'
$ws.Range("F20").Value = ' Restore Reach Function'
$ws.Range("F21").Value = ' Restore Reach Function'
$ws.Range("F22").Value = ' Restore Reach Function'
$ws.Range("F23").Value = ' Restore Reach Function'
$ws.Range("F24").Value = ' Restore Reach Function'
$ws.Range("F25").Value = ' Address Limiting Factors'
$ws.Range("F26").Value = ' Address Limiting Factors'
$ws.Range("F27").Value = ' Address Limiting Factors'
$ws.Range("F28").Value = ' Address Limiting Factors'
$ws.Range("F29").Value = ' Address Limiting Factors'
$ws.Range("F30").Value = ' Address Limiting Factors'
$ws.Range("F31").Value = ' Address Limiting Factors'
$ws.Range("F32").Value = ' Restore Reach Function'
$ws.Range("F33").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F34").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F35").Value = ' Restore Reach Function'
$ws.Range("F36").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F37").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F38").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F39").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F40").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F41").Value = ' Restore Reach Function'
$ws.Range("F42").Value = ' Address Limiting Factors'
$ws.Range("F43").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F44").Value = ' Restore Reach Function'
$ws.Range("F45").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F46").Value = ' Restore Reach Function'
$ws.Range("F47").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F48").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F49").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F50").Value = 'Restore Fish Passage'
$ws.Range("F51").Value = ' Restore Reach Function'
$ws.Range("F52").Value = 'Restore Fish Passage'
$ws.Range("F53").Value = ' Restore Reach Function'
$ws.Range("F54").Value = ' Restore Reach Function'
$ws.Range("F55").Value = ' Restore Reach Function'
$ws.Range("F56").Value = ' Restore Reach Function'
$ws.Range("F57").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F58").Value = ' Restore Reach Function'
$ws.Range("F59").Value = ' Restore Reach Function'
$ws.Range("F60").Value = ' Restore Reach Function'
$ws.Range("F61").Value = ' Restore Reach Function'
$ws.Range("F62").Value = ' Restore Reach Function'
$ws.Range("F63").Value = ' Restore Reach Function'
$ws.Range("F64").Value = ' Restore Reach Function'
$ws.Range("F65").Value = ' Restore Reach Function'
$ws.Range("F66").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F67").Value = ' Restore Reach Function'
$ws.Range("F68").Value = ' Restore Reach Function'
$ws.Range("F69").Value = ' Restore Reach Function'
$ws.Range("F70").Value = ' Restore Reach Function'
$ws.Range("F71").Value = ' Restore Reach Function'
$ws.Range("F72").Value = 'Restore Fish Passage'
$ws.Range("F73").Value = 'Restore Fish Passage'
$ws.Range("F74").Value = ' Address Limiting Factors'
$ws.Range("F75").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F76").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F77").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F78").Value = ' Address Limiting Factors'
$ws.Range("F79").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F80").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F81").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F82").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F83").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F84").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F85").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F86").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F87").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F88").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F89").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F90").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F91").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F92").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F93").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F94").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F95").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F96").Value = ' Restore Reach Function'
$ws.Range("F97").Value = ' Restore Reach Function'
$ws.Range("F98").Value = ' Restore Reach Function'
$ws.Range("F99").Value = ' Restore Reach Function'
$ws.Range("F100").Value = ' Restore Reach Function'
$ws.Range("F101").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F102").Value = ' Restore Reach Function'
$ws.Range("F103").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F104").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F105").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F106").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F107").Value = ' Address Limiting Factors'
$ws.Range("F108").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F109").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F110").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F111").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F112").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F113").Value = ' Address Limiting Factors'
$ws.Range("F114").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F115").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F116").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F117").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F118").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F119").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F120").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F121").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F122").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F123").Value = ' Restore Reach Function, Address Limiting Factors'
$ws.Range("F124").Value = ' Restore Reach Function'
$ws.Range("F125").Value = ' Restore Reach Function'
